# Refresh the crypto price/volume snapshot (values scraped from coinranking.com).
# D-column prices are forced back to text (leading apostrophe => quote-prefixed
# text entry) so values such as "325.55" or "0.4588" aren't auto-converted to
# numbers by Excel -- the source sheet stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.766.12'
$ws.Range("E2").Value = '  +2.69%  '

$ws.Range("D3").Value = '''1.875.60'
$ws.Range("E3").Value = '  +2.54%  '

$ws.Range("D5").Value = '''325.55'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("E6").Value = '  +0.38%  '

$ws.Range("D7").Value = '''0.4588'
$ws.Range("E7").Value = '  -0.62%  '

$ws.Range("D8").Value = '''0.3871'
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").Value = '''0.07854'
$ws.Range("E9").Value = '  +0.38%  '

$ws.Range("D10").Value = '''0.9935'
$ws.Range("E10").Value = '  +3.71%  '

$ws.Range("D11").Value = '''21.79'
$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("D12").Value = '''1.907.14'
$ws.Range("E12").Value = '  +3.58%  '

$ws.Range("E13").Value = '  +1.86%  '

$ws.Range("D14").Value = '''5.715'
$ws.Range("E14").Value = '  +0.93%  '

$ws.Range("D15").Value = '''0.06963'
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").Value = '''88.50'
$ws.Range("E16").Value = '  +0.36%  '

$ws.Range("E17").Value = '  +0.47%  '

$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("D19").Value = '''16.84'
$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").Value = '''28.773.02'
$ws.Range("E21").Value = '  +2.67%  '

$ws.Range("D22").Value = '''5.280'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = '''11.06'
$ws.Range("E23").Value = '  +0.96%  '

$ws.Range("D24").Value = '''2.139'
$ws.Range("E24").Value = '  +3.08%  '

$ws.Range("D25").Value = '''2.103.04'
$ws.Range("E25").Value = '  +2.15%  '

$ws.Range("D26").Value = '''153.22'
$ws.Range("E26").Value = '  -0.75%  '

$ws.Range("D27").Value = '''19.26'
$ws.Range("E27").Value = '  +0.78%  '

$ws.Range("D28").Value = '''5.799'
$ws.Range("E28").Value = '  +2.19%  '

$ws.Range("D30").Value = '''119.11'
$ws.Range("E30").Value = '  +0.74%  '

$ws.Range("D31").Value = '''0.09319'
$ws.Range("E31").Value = '  +1.17%  '

$ws.Range("D32").Value = '''0.9195'
$ws.Range("E32").Value = '  -1.69%  '

$ws.Range("E33").Value = '  +1.21%  '

$ws.Range("E34").Value = '  +1.98%  '

$ws.Range("E35").Value = '  +0.64%  '

$ws.Range("D36").Value = '''0.05768'
$ws.Range("E36").Value = '  -0.82%  '

$ws.Range("D37").Value = '''1.151'
$ws.Range("E37").Value = '  +1.52%  '

$ws.Range("D38").Value = '''0.02071'
$ws.Range("E38").Value = '  -2.22%  '

$ws.Range("D39").Value = '''7.704'
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").Value = '''0.5648'
$ws.Range("E40").Value = '  +1.39%  '

$ws.Range("D41").Value = '''0.1792'
$ws.Range("E41").Value = '  +2.11%  '

$ws.Range("D42").Value = '''9.926'
$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("D43").Value = '''0.07221'
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").Value = '''11.84'
$ws.Range("E44").Value = '  +2.30%  '

$ws.Range("D45").Value = '''0.5303'
$ws.Range("E45").Value = '  +0.99%  '

$ws.Range("D46").Value = '''2.167'
$ws.Range("E46").Value = '  +3.72%  '

$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("D48").Value = '''113.78'
$ws.Range("E48").Value = '  +0.75%  '

$ws.Range("D49").Value = '''1.828'
$ws.Range("E49").Value = '  +0.25%  '

$ws.Range("D50").Value = '''2.409'
$ws.Range("E50").Value = '  +3.89%  '

$ws.Range("E51").Value = '  +0.44%  '
